$wb = $excel.ActiveWorkbook

# ===== Sheet 1: US Totals =====
$ws = $wb.Worksheets.Item("US Totals")

# Step 1: Insert new row at 69, pushing everything below down by 1
$ws.Rows.Item(69).Insert()

# Step 2: Fill new row 69 with date data
$ws.Range("A69").Value = "'03-29-2020"
$ws.Range("B69").Value = 140886
$ws.Range("C69").Formula = "=(B69/B68) - 1"
$ws.Range("D69").Value = 2467
$ws.Range("E69").Formula = "=(D69/D68) - 1"
$ws.Range("F69").Formula = "=AVERAGE(E63:E69)"

# Step 3: Update the Average (7-Day) row (now at row 71) formulas
$ws.Range("C71").Formula = "=AVERAGE(C63:C69)"
$ws.Range("E71").Formula = "=AVERAGE(E63:E69)"

# Step 4: Delete the "Today" row (now at row 74, shifted from old 73)
$ws.Rows.Item(74).Delete()

# Step 5: Update labels and values for rows 74-79 (Day+1 .. Day+6)
$ws.Range("A74").Value = "Day +1"
$ws.Range("D74").Value = 3170.543530407662
$ws.Range("E74").Value = 0.284961491690312

$ws.Range("A75").Value = "Day +2"
$ws.Range("D75").Value = 4074.724879695938
$ws.Range("E75").Value = 0.6513984596676126

$ws.Range("A76").Value = "Day +3"
$ws.Range("D76").Value = 5236.762304625495
$ws.Range("E76").Value = 1.122415889744629

$ws.Range("A77").Value = "Day +4"
$ws.Range("D77").Value = 6730.191668104208
$ws.Range("E77").Value = 1.728009728415079

$ws.Range("A78").Value = "Day +5"
$ws.Range("D78").Value = 8649.519923677075
$ws.Range("E78").Value = 2.505877584110255

$ws.Range("A79").Value = "Day +6"
$ws.Range("D79").Value = 11116.20568915546
$ws.Range("E79").Value = 3.505877584110255

# Step 6: Add new row 80 "Day +7"
$ws.Range("A80").Value = "Day +7"
$ws.Range("D80").NumberFormat = "0"
$ws.Range("D80").Value = 14286.34537107122
$ws.Range("E80").NumberFormat = "00.00%"
$ws.Range("E80").Value = 4.790839075800568

# ===== Sheet 2: GA Totals =====
$ws = $wb.Worksheets.Item("GA Totals")

# Step 1: Insert new row at 69, pushing everything below down by 1
$ws.Rows.Item(69).Insert()

# Step 2: Fill new row 69 with date data
$ws.Range("A69").Value = "'03-29-2020"
$ws.Range("B69").Value = 2651
$ws.Range("C69").Formula = "=(B69/B68) - 1"
$ws.Range("D69").Value = 80
$ws.Range("E69").Formula = "=(D69/D68) - 1"
$ws.Range("F69").Formula = "=AVERAGE(E63:E69)"

# Step 3: Update the Average (7-Day) row (now at row 71) formulas
$ws.Range("C71").Formula = "=AVERAGE(C63:C69)"
$ws.Range("E71").Formula = "=AVERAGE(E63:E69)"

# Step 4: Delete the "Today" row (now at row 74, shifted from old 73)
$ws.Rows.Item(74).Delete()

# Step 5: Update labels and values for rows 74-79 (Day+1 .. Day+6)
$ws.Range("A74").Value = "Day +1"
$ws.Range("D74").Value = 94.86718426501037
$ws.Range("E74").Value = 0.175

$ws.Range("A75").Value = "Day +2"
$ws.Range("D75").Value = 112.4972831296429
$ws.Range("E75").Value = 0.3999999999999999

$ws.Range("A76").Value = "Day +3"
$ws.Range("D76").Value = 133.4037560996609
$ws.Range("E76").Value = 0.6625000000000001

$ws.Range("A77").Value = "Day +4"
$ws.Range("D77").Value = 158.1954838943879
$ws.Range("E77").Value = 0.9750000000000001

$ws.Range("A78").Value = "Day +5"
$ws.Range("D78").Value = 187.5945015062672
$ws.Range("E78").Value = 1.3375

$ws.Range("A79").Value = "Day +6"
$ws.Range("D79").Value = 222.4570267687227
$ws.Range("E79").Value = 1.775

# Step 6: Add new row 80 "Day +7"
$ws.Range("A80").Value = "Day +7"
$ws.Range("D80").NumberFormat = "0"
$ws.Range("D80").Value = 263.7983968689346
$ws.Range("E80").NumberFormat = "00.00%"
$ws.Range("E80").Value = 2.2875

# ===== Sheet 3: SC Totals =====
$ws = $wb.Worksheets.Item("SC Totals")

# Step 1: Insert new row at 69, pushing everything below down by 1
$ws.Rows.Item(69).Insert()

# Step 2: Fill new row 69 with date data
$ws.Range("A69").Value = "'03-29-2020"
$ws.Range("B69").Value = 774
$ws.Range("C69").Formula = "=(B69/B68) - 1"
$ws.Range("D69").Value = 16
$ws.Range("E69").Formula = "=(D69/D68) - 1"
$ws.Range("F69").Formula = "=AVERAGE(E63:E69)"

# Step 3: Update the Average (7-Day) row (now at row 71) formulas
$ws.Range("C71").Formula = "=AVERAGE(C63:C69)"
$ws.Range("E71").Formula = "=AVERAGE(E63:E69)"

# Step 4: Delete the "Today" row (now at row 74, shifted from old 73)
$ws.Rows.Item(74).Delete()

# Step 5: Update labels and values for rows 74-79 (Day+1 .. Day+6)
$ws.Range("A74").Value = "Day +1"
$ws.Range("D74").Value = 20.61105878248735
$ws.Range("E74").Value = 0.25

$ws.Range("A75").Value = "Day +2"
$ws.Range("D75").Value = 26.55098400844682
$ws.Range("E75").Value = 0.625

$ws.Range("A76").Value = "Day +3"
$ws.Range("D76").Value = 34.2027432581862
$ws.Range("E76").Value = 1.125

$ws.Range("A77").Value = "Day +4"
$ws.Range("D77").Value = 44.05967198854992
$ws.Range("E77").Value = 1.75

$ws.Range("A78").Value = "Day +5"
$ws.Range("D78").Value = 56.75728058081963
$ws.Range("E78").Value = 2.5

$ws.Range("A79").Value = "Day +6"
$ws.Range("D79").Value = 73.11422789908758
$ws.Range("E79").Value = 3.5625

# Step 6: Add new row 80 "Day +7"
$ws.Range("A80").Value = "Day +7"
$ws.Range("D80").NumberFormat = "0"
$ws.Range("D80").Value = 94.18510306651694
$ws.Range("E80").NumberFormat = "00.00%"
$ws.Range("E80").Value = 4.875

# ===== Sheet 4: NY Totals =====
$ws = $wb.Worksheets.Item("NY Totals")

# Step 1: Insert new row at 69, pushing everything below down by 1
$ws.Rows.Item(69).Insert()

# Step 2: Fill new row 69 with date data
$ws.Range("A69").Value = "'03-29-2020"
$ws.Range("B69").Value = 59648
$ws.Range("C69").Formula = "=(B69/B68) - 1"
$ws.Range("D69").Value = 965
$ws.Range("E69").Formula = "=(D69/D68) - 1"
$ws.Range("F69").Formula = "=AVERAGE(E63:E69)"

# Step 3: Update the Average (7-Day) row (now at row 71) formulas
$ws.Range("C71").Formula = "=AVERAGE(C63:C69)"
$ws.Range("E71").Formula = "=AVERAGE(E63:E69)"

# Step 4: Delete the "Today" row (now at row 74, shifted from old 73)
$ws.Rows.Item(74).Delete()

# Step 5: Update labels and values for rows 74-79 (Day+1 .. Day+6)
$ws.Range("A74").Value = "Day +1"
$ws.Range("D74").Value = 1304.589857772497
$ws.Range("E74").Value = 0.3512953367875649

$ws.Range("A75").Value = "Day +2"
$ws.Range("D75").Value = 1763.683623837164
$ws.Range("E75").Value = 0.8269430051813471

$ws.Range("A76").Value = "Day +3"
$ws.Range("D76").Value = 2384.335510857417
$ws.Range("E76").Value = 1.470466321243523

$ws.Range("A77").Value = "Day +4"
$ws.Range("D77").Value = 3223.398886001442
$ws.Range("E77").Value = 2.339896373056995

$ws.Range("A78").Value = "Day +5"
$ws.Range("D78").Value = 4357.734190914659
$ws.Range("E78").Value = 3.515025906735751

$ws.Range("A79").Value = "Day +6"
$ws.Range("D79").Value = 5891.249563042177
$ws.Range("E79").Value = 5.104663212435233

# Step 6: Add new row 80 "Day +7"
$ws.Range("A80").Value = "Day +7"
$ws.Range("D80").NumberFormat = "0"
$ws.Range("D80").Value = 7964.419097980806
$ws.Range("E80").NumberFormat = "00.00%"
$ws.Range("E80").Value = 7.252849740932643

# ===== Sheet 5: NC Totals =====
$ws = $wb.Worksheets.Item("NC Totals")

# Step 1: Insert new row at 69, pushing everything below down by 1
$ws.Rows.Item(69).Insert()

# Step 2: Fill new row 69 with date data
$ws.Range("A69").Value = "'03-29-2020"
$ws.Range("B69").Value = 1191
$ws.Range("C69").Formula = "=(B69/B68) - 1"
$ws.Range("D69").Value = 7
$ws.Range("E69").Formula = "=(D69/D68) - 1"
$ws.Range("F69").Formula = "=AVERAGE(E63:E69)"

# Step 3: Update the Average (7-Day) row (now at row 71) formulas
$ws.Range("C71").Formula = "=AVERAGE(C63:C69)"
$ws.Range("E71").Formula = "=AVERAGE(E63:E69)"

# Step 4: Delete the "Today" row (now at row 74, shifted from old 73)
$ws.Rows.Item(74).Delete()

# Step 5: Update labels and values for rows 74-79 (Day+1 .. Day+6)
$ws.Range("A74").Value = "Day +1"
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = -1

$ws.Range("A75").Value = "Day +2"
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = -1

$ws.Range("A76").Value = "Day +3"
$ws.Range("D76").Value = 0
$ws.Range("E76").Value = -1

$ws.Range("A77").Value = "Day +4"
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = -1

$ws.Range("A78").Value = "Day +5"
$ws.Range("D78").Value = 0
$ws.Range("E78").Value = -1

$ws.Range("A79").Value = "Day +6"
$ws.Range("D79").Value = 0
$ws.Range("E79").Value = -1

# Step 6: Add new row 80 "Day +7"
$ws.Range("A80").Value = "Day +7"
$ws.Range("D80").NumberFormat = "0"
$ws.Range("D80").Value = 0
$ws.Range("E80").NumberFormat = "00.00%"
$ws.Range("E80").Value = -1

# ===== Sheet 6: WA Totals =====
$ws = $wb.Worksheets.Item("WA Totals")

# Step 1: Insert new row at 69, pushing everything below down by 1
$ws.Rows.Item(69).Insert()

# Step 2: Fill new row 69 with date data
$ws.Range("A69").Value = "'03-29-2020"
$ws.Range("B69").Value = 4465
$ws.Range("C69").Formula = "=(B69/B68) - 1"
$ws.Range("D69").Value = 198
$ws.Range("E69").Formula = "=(D69/D68) - 1"
$ws.Range("F69").Formula = "=AVERAGE(E63:E69)"

# Step 3: Update the Average (7-Day) row (now at row 71) formulas
$ws.Range("C71").Formula = "=AVERAGE(C63:C69)"
$ws.Range("E71").Formula = "=AVERAGE(E63:E69)"

# Step 4: Delete the "Today" row (now at row 74, shifted from old 73)
$ws.Rows.Item(74).Delete()

# Step 5: Update labels and values for rows 74-79 (Day+1 .. Day+6)
$ws.Range("A74").Value = "Day +1"
$ws.Range("D74").Value = 219.4862014345704
$ws.Range("E74").Value = 0.106060606060606

$ws.Range("A75").Value = "Day +2"
$ws.Range("D75").Value = 243.3040031322061
$ws.Range("E75").Value = 0.2272727272727273

$ws.Range("A76").Value = "Day +3"
$ws.Range("D76").Value = 269.7064214207714
$ws.Range("E76").Value = 0.3585858585858586

$ws.Range("A77").Value = "Day +4"
$ws.Range("D77").Value = 298.9739289906898
$ws.Range("E77").Value = 0.505050505050505

$ws.Range("A78").Value = "Day +5"
$ws.Range("D78").Value = 331.4174343542196
$ws.Range("E78").Value = 0.6717171717171717

$ws.Range("A79").Value = "Day +6"
$ws.Range("D79").Value = 367.3815846242362
$ws.Range("E79").Value = 0.8535353535353536

# Step 6: Add new row 80 "Day +7"
$ws.Range("A80").Value = "Day +7"
$ws.Range("D80").NumberFormat = "0"
$ws.Range("D80").Value = 407.248426697913
$ws.Range("E80").NumberFormat = "00.00%"
$ws.Range("E80").Value = 1.055555555555555

# ===== Sheet 7: FL Totals =====
$ws = $wb.Worksheets.Item("FL Totals")

# Step 1: Insert new row at 69, pushing everything below down by 1
$ws.Rows.Item(69).Insert()

# Step 2: Fill new row 69 with date data
$ws.Range("A69").Value = "'03-29-2020"
$ws.Range("B69").Value = 4246
$ws.Range("C69").Formula = "=(B69/B68) - 1"
$ws.Range("D69").Value = 56
$ws.Range("E69").Formula = "=(D69/D68) - 1"
$ws.Range("F69").Formula = "=AVERAGE(E63:E69)"

# Step 3: Update the Average (7-Day) row (now at row 71) formulas
$ws.Range("C71").Formula = "=AVERAGE(C63:C69)"
$ws.Range("E71").Formula = "=AVERAGE(E63:E69)"

# Step 4: Delete the "Today" row (now at row 74, shifted from old 73)
$ws.Rows.Item(74).Delete()

# Step 5: Update labels and values for rows 74-79 (Day+1 .. Day+6)
$ws.Range("A74").Value = "Day +1"
$ws.Range("D74").Value = 69.68042767383098
$ws.Range("E74").Value = 0.2321428571428572

$ws.Range("A75").Value = "Day +2"
$ws.Range("D75").Value = 86.70289287157124
$ws.Range("E75").Value = 0.5357142857142858

$ws.Range("A76").Value = "Day +3"
$ws.Range("D76").Value = 107.8838331401685
$ws.Range("E76").Value = 0.9107142857142858

$ws.Range("A77").Value = "Day +4"
$ws.Range("D77").Value = 134.2391362910564
$ws.Range("E77").Value = 1.392857142857143

$ws.Range("A78").Value = "Day +5"
$ws.Range("D78").Value = 167.0328647736874
$ws.Range("E78").Value = 1.982142857142857

$ws.Range("A79").Value = "Day +6"
$ws.Range("D79").Value = 207.8378830895663
$ws.Range("E79").Value = 2.696428571428572

# Step 6: Add new row 80 "Day +7"
$ws.Range("A80").Value = "Day +7"
$ws.Range("D80").NumberFormat = "0"
$ws.Range("D80").Value = 258.6112960804404
$ws.Range("E80").NumberFormat = "00.00%"
$ws.Range("E80").Value = 3.607142857142857

# ===== Sheet 8: CA Totals =====
$ws = $wb.Worksheets.Item("CA Totals")

# Step 1: Insert new row at 69, pushing everything below down by 1
$ws.Rows.Item(69).Insert()

# Step 2: Fill new row 69 with date data
$ws.Range("A69").Value = "'03-29-2020"
$ws.Range("B69").Value = 5852
$ws.Range("C69").Formula = "=(B69/B68) - 1"
$ws.Range("D69").Value = 124
$ws.Range("E69").Formula = "=(D69/D68) - 1"
$ws.Range("F69").Formula = "=AVERAGE(E63:E69)"

# Step 3: Update the Average (7-Day) row (now at row 71) formulas
$ws.Range("C71").Formula = "=AVERAGE(C63:C69)"
$ws.Range("E71").Formula = "=AVERAGE(E63:E69)"

# Step 4: Delete the "Today" row (now at row 74, shifted from old 73)
$ws.Rows.Item(74).Delete()

# Step 5: Update labels and values for rows 74-79 (Day+1 .. Day+6)
$ws.Range("A74").Value = "Day +1"
$ws.Range("D74").Value = 152.0981245179827
$ws.Range("E74").Value = 0.2258064516129032

$ws.Range("A75").Value = "Day +2"
$ws.Range("D75").Value = 186.5632216281271
$ws.Range("E75").Value = 0.5

$ws.Range("A76").Value = "Day +3"
$ws.Range("D76").Value = 228.8380331747651
$ws.Range("E76").Value = 0.8387096774193548

$ws.Range("A77").Value = "Day +4"
$ws.Range("D77").Value = 280.6922230989167
$ws.Range("E77").Value = 1.258064516129032

$ws.Range("A78").Value = "Day +5"
$ws.Range("D78").Value = 344.2964572591
$ws.Range("E78").Value = 1.774193548387097

$ws.Range("A79").Value = "Day +6"
$ws.Range("D79").Value = 422.3132695749587
$ws.Range("E79").Value = 2.403225806451613

# Step 6: Add new row 80 "Day +7"
$ws.Range("A80").Value = "Day +7"
$ws.Range("D80").NumberFormat = "0"
$ws.Range("D80").Value = 518.0085182371649
$ws.Range("E80").NumberFormat = "00.00%"
$ws.Range("E80").Value = 3.17741935483871
